# 612 report work distribution.xlsx
# Commit: "Added work to convert sklearn model to pyspark machine learning"
#
# Adds a new row-5 task entry in column E ("Final Project Work" list)
# describing the new work item, and moves the active selection down to
# E6 (where the author's cursor ended up after typing the new entry and
# pressing Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string entry: the work item that was added.
$ws.Range("E5").Value = "Convert model from sklearn to pyspark"

# After entering the value in E5, the author's selection moved to E6.
$ws.Range("E6").Select()
